$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Limited data set" -------------------------------------------
# Alphabetically it belongs between "Instrument" (row 27) and
# "Longitudinal" (row 28), so push everything at/after row 28 down by one.
$ws.Rows.Item(28).Insert()
$ws.Range("C28").Value = "A dataset in which 15 of the 18 HIPAA protected identifiers have been removed. Age, dates, and city/state/zipcode can remain."
$ws.Range("A28").Value = "Limited data set"

# --- Insert "Restricted-use data" ----------------------------------------
# Alphabetically it belongs between "Research" (now row 46) and
# "Scale" (now row 47), so push everything at/after row 47 down by one.
$ws.Rows.Item(47).Insert()
$ws.Range("C47").Value = "A dataset that contains sensitive, confidential, or proprietary information."
$ws.Range("A47").Value = "Restricted-use data"

# --- Re-apply the existing alphabetical sort over the whole table --------
# (keeps the worksheet's recorded sort range/condition in sync with the two
# newly inserted rows; the data is already in order so this is a no-op on
# the values themselves)
$sortRange = $ws.Range("A2:C60")
$keyRange = $ws.Range("A2:A60")
$s = $ws.Sort
$s.SortFields.Clear()
$s.SortFields.Add($keyRange) | Out-Null
$s.SetRange($sortRange)
$s.Apply()

# --- Leave the view/selection where the author left it -------------------
$ws.Activate()
$ws.Range("C34").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
